# edit.ps1 - applies the "dark matter" -> "chemistry" rewrite described in the diff.
$d = $word.ActiveDocument
$br = [char]11

# ---------------------------------------------------------------------------
# 1) Title
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Unraveling the Enigma of Dark Matter", $true, $false, $false, $false, $false, $true, 1, $false, "Navigating the Marvelous World of Chemistry: Unlocking the Secrets of Matter", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2) Author name: "Alexander Weiss" -> "Dr. Sarah Montgomery"
# ---------------------------------------------------------------------------
$authorPara = $d.Paragraphs(2)
$authorRange = $d.Range($authorPara.Range.Start, $authorPara.Range.End - 1)
$authorRange.Text = "Dr" + "." + " Sarah Montgomery"

# ---------------------------------------------------------------------------
# 3) Email: "alexanderweiss@protonmail.com" -> "sarah.montgomery@educationalhaven.org"
# ---------------------------------------------------------------------------
$emailPara = $d.Paragraphs(3)
$emailRange = $d.Range($emailPara.Range.Start, $emailPara.Range.End - 1)
$emailRange.Text = "sarah" + "." + "montgomery@educationalhaven" + "." + "org"

# ---------------------------------------------------------------------------
# 4) Body paragraph (long rewrite, with manual line breaks)
# ---------------------------------------------------------------------------
$bodyPara = $d.Paragraphs(5)
$bodyRange = $d.Range($bodyPara.Range.Start, $bodyPara.Range.End - 1)
$bodyRange.Text = "In the vast expanse of scientific disciplines, chemistry stands as a beacon of discovery, revealing the fundamental principles that govern the composition and behavior of matter. It is a symphony of elements, a dance of molecules, a journey through the intricate tapestry of chemical reactions. As we delve into this captivating realm, we embark on an expedition to unravel the secrets of the universe's most fundamental building blocks." + $br + "" + $br + "From the tiniest atoms to the grandest molecules, chemistry holds the key to understanding the world around us. Its principles are woven into the fabric of our existence, shaping the properties of everything from the air we breathe to the food we consume. Chemistry empowers us to create life-saving medicines, harness energy from renewable sources, and develop innovative materials that drive technological advancements." + $br + "" + $br + "The study of chemistry is an odyssey into the unknown, a voyage where each experiment is a step towards unraveling the mysteries of matter. It challenges us to think critically, to question the world around us, and to seek answers through meticulous observation and experimentation. Chemistry is a gateway to a world of wonder, a testament to the boundless curiosity of the human spirit." + $br + "" + $br + "Body:" + $br + "" + $br + "In the realm of chemistry, we dissect the structure of substances, discovering the elements that compose them and the forces that hold them together. We explore the intricate dance of chemical reactions, witnessing the transformation of one substance into another. We unravel the secrets of chemical bonding, the language of molecules, deciphering the patterns and properties that determine the behavior of matter." + $br + "" + $br + "Through the lens of chemistry, we gain insights into the natural world, comprehending the intricate processes that occur in ecosystems, the delicate balance of life, and the interdependence of all living things. We investigate the chemistry of life, unraveling the secrets of cellular respiration, photosynthesis, and DNA replication. We delve into the mysteries of human health, exploring the intricate workings of our bodies and the chemical reactions that keep us alive." + $br + "" + $br + "Chemistry also has a profound impact on our daily lives, influencing everything from the food we eat to the medicines we take. We harness the power of chemistry to create innovative materials, from lightweight alloys to durable plastics, revolutionizing industries and transforming our world. Chemistry enables us to develop sustainable energy sources, reduce environmental pollution, and address global challenges such as climate change."

# ---------------------------------------------------------------------------
# 5) Summary paragraph text
# ---------------------------------------------------------------------------
$summaryPara = $d.Paragraphs(7)
$summaryRange = $d.Range($summaryPara.Range.Start, $summaryPara.Range.End - 1)
$summaryRange.Text = "Chemistry is a captivating and transformative discipline that unlocks the secrets of matter and its interactions. It empowers us to understand the world around us, from the tiniest atoms to the grandest molecules. Chemistry has revolutionized our lives, leading to life-saving medicines, innovative materials, and sustainable technologies. It challenges us to think critically, question the world around us, and seek answers through experimentation. Chemistry is a testament to the boundless curiosity of the human spirit, a journey of discovery that continues to reshape our understanding of the universe."

# ---------------------------------------------------------------------------
# 6) Append a new empty paragraph at the very end of the document
# ---------------------------------------------------------------------------
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$lastPara.Range.InsertParagraphAfter()

# ---------------------------------------------------------------------------
# 7) Global font fix: "TimesNewToman" -> "Times New Roman"
# ---------------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    try {
        $p.Range.Font.Name = "Times New Roman"
    } catch {
        # paragraphs with no runs (e.g. a bare <w:p/>) have nothing to restyle
    }
}

Write-Host "Edit complete"
